$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.439.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "'1.917.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").Value = "'324.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "'0.4822"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "'0.4078"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "'0.08240"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").Value = "'1.016"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "'23.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "'1.915.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "'6.083"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").Value = "'91.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "'0.06815"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "'17.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "'29.474.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'5.657"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "'11.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").Value = "'2.176"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "'2.140.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "'6.654"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.16%  "
$ws.Range("D27").Value = "'155.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "'20.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").Value = "'2.115"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").Value = "'120.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("D31").Value = "'1.021"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'0.09599"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("D33").Value = "'5.698"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.00%  "
$ws.Range("D34").Value = "'3.552"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").Value = "'1.373"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").Value = "'0.02289"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "'0.06111"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "'1.184"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'8.090"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").Value = "'0.5988"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").Value = "'10.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.07%  "
$ws.Range("D42").Value = "'0.1850"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.280"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.402"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'12.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.88%  "
$ws.Range("D46").Value = "'0.07597"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "'0.5592"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "'1.957"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "'118.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.84%  "
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("D51").Value = "'72.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
